# Scheduled-runner data refresh for the Ultima_Profits leve-crafting
# workbook: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) on rows whose market data changed, across every job
# sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1714.4762
$ws.Range("I40").Value = 1607.9231
$ws.Range("J40").Value = 1887.625
$ws.Range("K40").Value = 1607.9231
$ws.Range("L40").Value = 1887.625
$ws.Range("M40").Value = -1432.9231
$ws.Range("N40").Value = -2237.625

$ws.Range("H129").Value = 2786.3948
$ws.Range("I129").Value = 443.9091
$ws.Range("J129").Value = 3740.7407
$ws.Range("K129").Value = 1331.7273
$ws.Range("L129").Value = 11222.2221
$ws.Range("M129").Value = 3668.2727
$ws.Range("N129").Value = -21222.2221

$ws.Range("H138").Value = 5436848
$ws.Range("J138").Value = 9617865
$ws.Range("L138").Value = 28853595
$ws.Range("N138").Value = -28863875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10044.333
$ws.Range("I32").Value = 10750.659
$ws.Range("K32").Value = 10750.659
$ws.Range("M32").Value = -10463.659

$ws.Range("H61").Value = 13516267
$ws.Range("I61").Value = 17243892
$ws.Range("J61").Value = 3626.75
$ws.Range("K61").Value = 17243892
$ws.Range("L61").Value = 3626.75
$ws.Range("M61").Value = -17243680
$ws.Range("N61").Value = -4050.75

$ws.Range("H136").Value = 13516267
$ws.Range("I136").Value = 17243892
$ws.Range("J136").Value = 3626.75
$ws.Range("K136").Value = 51731676
$ws.Range("L136").Value = 10880.25
$ws.Range("M136").Value = -51729126
$ws.Range("N136").Value = -15980.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1135.1666
$ws.Range("I99").Value = 1150
$ws.Range("J99").Value = 1127.75
$ws.Range("K99").Value = 1150
$ws.Range("L99").Value = 1127.75
$ws.Range("M99").Value = 348
$ws.Range("N99").Value = -4123.75

$ws.Range("H105").Value = 4353.147
$ws.Range("I105").Value = 3200.9
$ws.Range("J105").Value = 4833.25
$ws.Range("K105").Value = 3200.9
$ws.Range("L105").Value = 4833.25
$ws.Range("M105").Value = -1453.9
$ws.Range("N105").Value = -8327.25

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1063.7368
$ws.Range("I16").Value = 882.3
$ws.Range("J16").Value = 1265.3334
$ws.Range("K16").Value = 882.3
$ws.Range("L16").Value = 1265.3334
$ws.Range("M16").Value = -595.3
$ws.Range("N16").Value = -1839.3334

$ws.Range("H31").Value = 6670231.5
$ws.Range("I31").Value = 3609.5334
$ws.Range("J31").Value = 66669830
$ws.Range("K31").Value = 3609.5334
$ws.Range("L31").Value = 66669830
$ws.Range("M31").Value = -3314.5334
$ws.Range("N31").Value = -66670420

$ws.Range("H34").Value = 6670231.5
$ws.Range("I34").Value = 3609.5334
$ws.Range("J34").Value = 66669830
$ws.Range("K34").Value = 3609.5334
$ws.Range("L34").Value = 66669830
$ws.Range("M34").Value = -3407.5334
$ws.Range("N34").Value = -66670234

$ws.Range("H99").Value = 1216.5
$ws.Range("I99").Value = 1064.2727
$ws.Range("J99").Value = 1368.7273
$ws.Range("K99").Value = 1064.2727
$ws.Range("L99").Value = 1368.7273
$ws.Range("M99").Value = 433.7273
$ws.Range("N99").Value = -4364.7273

$ws.Range("H107").Value = 543.3333
$ws.Range("I107").Value = 463.36365
$ws.Range("J107").Value = 669
$ws.Range("K107").Value = 463.36365
$ws.Range("L107").Value = 669
$ws.Range("M107").Value = 1456.63635
$ws.Range("N107").Value = -4509

$ws.Range("H113").Value = 1063.7368
$ws.Range("I113").Value = 882.3
$ws.Range("J113").Value = 1265.3334
$ws.Range("K113").Value = 882.3
$ws.Range("L113").Value = 1265.3334
$ws.Range("M113").Value = 1287.7
$ws.Range("N113").Value = -5605.3334

$ws.Range("H126").Value = 1216.5
$ws.Range("I126").Value = 1064.2727
$ws.Range("J126").Value = 1368.7273
$ws.Range("K126").Value = 3192.8181
$ws.Range("L126").Value = 4106.1819
$ws.Range("M126").Value = -722.8181
$ws.Range("N126").Value = -9046.1819

$ws.Range("H132").Value = 3318.5386
$ws.Range("I132").Value = 2650.5334
$ws.Range("K132").Value = 7951.600199999999
$ws.Range("M132").Value = -5421.600199999999

$ws.Range("H134").Value = 852475.3
$ws.Range("I134").Value = 2191.6924
$ws.Range("K134").Value = 6575.0772
$ws.Range("M134").Value = -4040.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1700.4667
$ws.Range("J22").Value = 2300.7
$ws.Range("L22").Value = 6902.099999999999
$ws.Range("N22").Value = -7240.099999999999

$ws.Range("H27").Value = 1700.4667
$ws.Range("J27").Value = 2300.7
$ws.Range("L27").Value = 6902.099999999999
$ws.Range("N27").Value = -7106.099999999999

$ws.Range("H33").Value = 106.125
$ws.Range("J33").Value = 100
$ws.Range("L33").Value = 600
$ws.Range("N33").Value = -1166

$ws.Range("H57").Value = 2300
$ws.Range("I57").Value = 500
$ws.Range("J57").Value = 3200
$ws.Range("K57").Value = 1500
$ws.Range("L57").Value = 9600
$ws.Range("M57").Value = -941
$ws.Range("N57").Value = -10718

$ws.Range("H58").Value = 2155.5557
$ws.Range("I58").Value = 466.66666
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1399.99998
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -1271.99998
$ws.Range("N58").Value = -9256

$ws.Range("H59").Value = 2385.7144
$ws.Range("I59").Value = 700
$ws.Range("J59").Value = 2666.6667
$ws.Range("K59").Value = 2100
$ws.Range("L59").Value = 8000.000100000001
$ws.Range("M59").Value = -1560
$ws.Range("N59").Value = -9080.000100000001

$ws.Range("H60").Value = 672.7778
$ws.Range("I60").Value = 594.375
$ws.Range("K60").Value = 1783.125
$ws.Range("M60").Value = -1532.125

$ws.Range("H61").Value = 579
$ws.Range("I61").Value = 334.66666
$ws.Range("J61").Value = 701.1667
$ws.Range("K61").Value = 1003.99998
$ws.Range("L61").Value = 2103.5001
$ws.Range("M61").Value = -788.9999799999999
$ws.Range("N61").Value = -2533.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 9985.5
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 9985.5
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 9985.5
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -10545.5

$ws.Range("H50").Value = 9985.5
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 9985.5
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 9985.5
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -10981.5

$ws.Range("H122").Value = 6669516.5
$ws.Range("I122").Value = 11113194
$ws.Range("J122").Value = 4001
$ws.Range("K122").Value = 33339582
$ws.Range("L122").Value = 12003
$ws.Range("M122").Value = -33337132
$ws.Range("N122").Value = -16903

$ws.Range("H138").Value = 56728
$ws.Range("J138").Value = 56728
$ws.Range("L138").Value = 56728
$ws.Range("N138").Value = -67008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1272.2727
$ws.Range("I61").Value = 1110.5555
$ws.Range("K61").Value = 1110.5555
$ws.Range("M61").Value = -908.5554999999999

$ws.Range("H113").Value = 1272.2727
$ws.Range("I113").Value = 1110.5555
$ws.Range("K113").Value = 1110.5555
$ws.Range("M113").Value = 1059.4445

$ws.Range("H139").Value = 39569.4
$ws.Range("J139").Value = 39449.332
$ws.Range("L139").Value = 39449.332
$ws.Range("N139").Value = -49729.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1344.5834
$ws.Range("I132").Value = 792.3333
$ws.Range("J132").Value = 3001.3333
$ws.Range("K132").Value = 2376.9999
$ws.Range("L132").Value = 9003.999899999999
$ws.Range("M132").Value = 153.0001000000002
$ws.Range("N132").Value = -14063.9999

$ws.Range("H136").Value = 1502.1305
$ws.Range("I136").Value = 1484.9412
$ws.Range("J136").Value = 1550.8334
$ws.Range("K136").Value = 4454.8236
$ws.Range("L136").Value = 4652.5002
$ws.Range("M136").Value = -1904.8236
$ws.Range("N136").Value = -9752.5002
